# "Generate Report for Handback" — the localization-status report is
# refreshed after a successful handback: the status flips from
# "Ready for handoff" to "Handed back: in sync with en-US", the handback
# timestamps for zh-cn and de-de move forward, and the stale-handback
# warning in the Error Detail column is cleared now that the handback is
# up to date with en-US.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# This is a shared string reused by Overview!E2/F2 and the per-language
# Status columns, so updating one propagates to all of them.
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"

# --- Latest Handback DateTime refreshed for each language ---
$zhcn.Range("K2").Value = "2016-08-13 07:02:05"
$dede.Range("K2").Value = "2016-08-13 07:02:15"

# --- Error Detail cleared now that the handback version is current ---
$zhcn.Range("P2").Value = ""
$dede.Range("P2").Value = ""

# --- Column widths re-fit to the refreshed content ---
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333332
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333332
